# Generate Report for Handback
# Updates the handback-status workbook with newly generated timestamps
# (and refreshes the zh-cn priority from "ht" to "mt").

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 41612351-...md row refreshed to the newly generated timestamp.
$wsOverview.Range("G2").Value = "2016-08-12 22:19:57"
$wsOverview.Range("G5").Value = "2016-08-12 22:19:57"

# Priority changed from "ht" to "mt" (shared by both the zh-cn and
# de-de sheets, since they referenced the same shared-string value).
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn sheet: handoff / handback datetimes for the 41612351-...md row
# bumped forward.
$wsZhCn.Range("H2").Value = "2016-08-12 22:19:50"
$wsZhCn.Range("H5").Value = "2016-08-12 22:19:50"
$wsZhCn.Range("K2").Value = "2016-08-12 22:20:18"
$wsZhCn.Range("K5").Value = "2016-08-12 22:20:18"

# de-de sheet: handoff datetime mirrors the Overview generate date, and
# handback datetime bumped forward as well.
$wsDeDe.Range("H2").Value = "2016-08-12 22:19:57"
$wsDeDe.Range("H5").Value = "2016-08-12 22:19:57"
$wsDeDe.Range("K2").Value = "2016-08-12 22:20:27"
$wsDeDe.Range("K5").Value = "2016-08-12 22:20:27"
